# Apply a cyclic re-shuffle of the observation rows 2, 3, 6, 8, 9.
# Content (Id, Taxonsorteringsordning, Rödlistade, TaxonId, Artnamn,
# Vetenskapligt namn, Auktor, Ost, Nord, Publik kommentar) moves along the
# cycle: row6 -> row2 -> row3 -> row9 -> row8 -> row6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in the change.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "AC")

# Snapshot current ("before") values for the affected rows, including the
# optional "AC" (Publik kommentar) cell which may be blank.
function Get-RowSnapshot([int]$row) {
    $snap = @{}
    foreach ($col in $cols) {
        $snap[$col] = $ws.Range("$col$row").Value2
    }
    return $snap
}

$row2 = Get-RowSnapshot 2
$row3 = Get-RowSnapshot 3
$row6 = Get-RowSnapshot 6
$row8 = Get-RowSnapshot 8
$row9 = Get-RowSnapshot 9

function Set-RowValues([int]$row, $values) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $values[$col]
    }
}

# Apply the cycle using the snapshots captured above so that every
# destination receives the correct source data.
Set-RowValues 2 $row6
Set-RowValues 3 $row2
Set-RowValues 6 $row8
Set-RowValues 8 $row9
Set-RowValues 9 $row3
